$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lead time")

# The staging table in rows 14:18 already holds the actual current-sale
# ("venta actual") figures for the "Venta Directa" offices. Promote those
# values into the placeholder rows 2:6 (which only carried office labels),
# matching the office order already used below (Africa, Agro Sudamerica,
# Agrosuper Asia, Agrosuper Brasil, Exportacion Directa).

$ws.Range("B2").Formula = "Africa"
$ws.Range("C2").Formula = "3.0000000000000004"
$ws.Range("D2").Formula = "7"
$ws.Range("E2:G2").ClearContents()
$ws.Range("H2").Formula = "=SUM(C2:D2)"

$ws.Range("B3").Formula = "Agro Sudamerica"
$ws.Range("C3").Formula = "3.0034758973294529"
$ws.Range("D3").Formula = "7"
$ws.Range("E3:G3").ClearContents()
$ws.Range("H3").Formula = "=SUM(C3:D3)"

$ws.Range("B4").Formula = "Agrosuper Asia"
$ws.Range("C4").Formula = "4.657346307107681"
$ws.Range("D4").Formula = "7"
$ws.Range("E4:G4").ClearContents()
$ws.Range("H4").Formula = "=SUM(C4:D4)"

$ws.Range("B5").Formula = "Agrosuper Brasil"
$ws.Range("C5").Formula = "4.2658490213417863"
$ws.Range("D5").ClearContents()
$ws.Range("E5:G5").ClearContents()
$ws.Range("H5").Formula = "=SUM(C5:D5)"

$ws.Range("B6").Formula = "Exportacion Directa"
$ws.Range("C6").Formula = "5.0000000000000027"
$ws.Range("D6").ClearContents()
$ws.Range("E6:G6").ClearContents()
$ws.Range("H6").Formula = "=SUM(C6:D6)"

# Row 12 used to carry a spare "Africa" office label with no values; that
# slot now just hosts the "Semanas" note that used to live on row 13.
$ws.Range("A12").ClearContents()
$ws.Range("B12").ClearContents()
$ws.Range("C12").Formula = "Semanas"

# The now-redundant staging rows (the old "Semanas" note row plus the
# 5-row table we just promoted above) are removed entirely.
$ws.Rows("13:18").Delete()
